# Adds "elderly" / "non_elderly" breakdown to the rates table.
# 1) Inserts a new (blank) column G, shifting old G/H/I -> H/I/J.
# 2) Resizes a few columns to the new widths.
# 3) Appends six new data rows (17-22) for the elderly / non-elderly groups.
# 4) Restores the selection to the cell the author ended up on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a blank column before the old "AAPC" column (G) -----------
# This shifts G->H, H->I, I->J and carries their number formats/values with
# them, so H/I/J end up with exactly the old G/H/I widths and cell styles.
$ws.Columns("G").Insert()

# The insert stamps a header-styled (but empty) placeholder into G1; the
# target layout has no cell at all there (only G2:G16 keep an empty,
# data-styled placeholder), so drop it explicitly.
$ws.Range("G1").Clear()

# --- 2) Column width tweaks ------------------------------------------------
# (H/I/J already carry the correct widths - 6.43 / 6.87 / 7.75 - from the
# column insert above, so only B/C/D/E need to be touched here.)
$ws.Columns("B").ColumnWidth = 43.64
$ws.Columns("C").ColumnWidth = 6.5
$ws.Range("D1:E1").EntireColumn.ColumnWidth = 5.33

# --- 3) Append the six new "elderly" / "non_elderly" rows ------------------
# Clone formatting from the last existing data row (16) so the new rows get
# the same per-column styles (integer / text / #,##0 / #,##0.00 + blank G).
$ws.Range("A16:J16").Copy()
$ws.Range("A17:J22").PasteSpecial(-4122)

$rows = @(
    @{ A = 15; B = "admissoes_gerais_non_elderly_tx_adjusted";         C = 4284; D = 3752; E = 3915; F = -8.61344537815126;  H = -4.50792455862232; I = -4.59019904070693; J = -4.42557912900052 },
    @{ A = 16; B = "admissoes_gerais_uti_non_elderly_tx_adjusted";     C = 149;  D = 170;  E = 248;  F = 66.4429530201342;   H = 30.4785145179664;  I = 29.9539815971385;  J = 31.0051646096701 },
    @{ A = 17; B = "admissoes_gerais_non_uti_non_elderly_tx_adjusted"; C = 4136; D = 3582; E = 3667; F = -11.3394584139265;  H = -5.96816560894209; I = -6.05120545564322; J = -5.88505236463991 },
    @{ A = 18; B = "admissoes_gerais_elderly_tx_adjusted";             C = 8553; D = 7268; E = 7624; F = -10.8616859581433;  H = -5.65466370282903; I = -5.90105148388402; J = -5.40763078237971 },
    @{ A = 19; B = "admissoes_gerais_uti_elderly_tx_adjusted";         C = 916;  D = 1068; E = 1228; F = 34.061135371179;    H = 15.7205950965359;  I = 14.9011753258851;  J = 16.5458585738211 },
    @{ A = 20; B = "admissoes_gerais_non_uti_elderly_tx_adjusted";     C = 7637; D = 6200; E = 6396; F = -16.2498363231636;  H = -8.66792163732836; I = -8.92490304718629; J = -8.41021511767397 }
)

$r = 17
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $r = $r + 1
}

# --- 4) Leave the selection where the author left it ------------------------
$ws.Range("P12").Select()
